# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# Mapping of worksheet name -> row number -> new value for column F
$updates = @{
    "展览" = @{
        2  = 64
        3  = 648
        4  = 218
        6  = 9858
        7  = 891
        10 = 3925
        11 = 173
        13 = 44
        16 = 548
        18 = 268
        19 = 1455
    }
    "全部类型" = @{
        2  = 64
        4  = 648
        5  = 218
        7  = 9858
        8  = 891
        11 = 3925
        12 = 173
        14 = 44
        17 = 548
        19 = 268
        20 = 1455
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsMap[$row]
    }
}
